# Apply updated crypto price/volume data (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '29.905.58'
$ws.Range('E2').Value = '  -0.14%  '

# Row 3
$ws.Range('D3').Value = '1.874.83'
$ws.Range('E3').Value = '  -0.98%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9984'
$ws.Range('E4').Value = '  -0.23%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7400'
$ws.Range('E5').Value = '  -4.74%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '242.39'
$ws.Range('E6').Value = '  -0.64%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.9991'

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3161'
$ws.Range('E8').Value = '  +0.99%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07198'
$ws.Range('E9').Value = '  -0.72%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '24.75'
$ws.Range('E10').Value = '  -4.26%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08363'
$ws.Range('E11').Value = '  -3.86%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.7507'
$ws.Range('E12').Value = '  -3.04%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.430'
$ws.Range('E13').Value = '  +0.30%  '

# Row 14
$ws.Range('D14').Value = '1.947.61'
$ws.Range('E14').Value = '  -5.47%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '92.62'
$ws.Range('E15').Value = '  -1.96%  '

# Row 16
$ws.Range('D16').Value = '29.913.36'
$ws.Range('E16').Value = '  -0.62%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '6.080'
$ws.Range('E17').Value = '  -1.88%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '13.58'
$ws.Range('E18').Value = '  -2.33%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '244.96'
$ws.Range('E19').Value = '  -0.32%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.000007831'
$ws.Range('E20').Value = '  -0.41%  '

# Row 21
$ws.Range('E21').Value = '  -0.28%  '

# Row 22
$ws.Range('D22').Value = '2.124.48'
$ws.Range('E22').Value = '  -6.90%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.995'
$ws.Range('E23').Value = '  -2.09%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.9984'
$ws.Range('E24').Value = '  -0.27%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1556'
$ws.Range('E25').Value = '  -5.27%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.272'
$ws.Range('E26').Value = '  -2.36%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '164.94'
$ws.Range('E27').Value = '  +0.96%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.62'
$ws.Range('E28').Value = '  -1.14%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.036'
$ws.Range('E29').Value = '  -0.95%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.515'
$ws.Range('E30').Value = '  +5.58%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.593'
$ws.Range('E31').Value = '  +1.60%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.534'
$ws.Range('E32').Value = '  -0.60%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.283'
$ws.Range('E33').Value = '  +3.72%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.05328'
$ws.Range('E34').Value = '  -2.72%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.238'
$ws.Range('E35').Value = '  -0.66%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7536'
$ws.Range('E36').Value = '  -0.19%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.002'
$ws.Range('E37').Value = '  -0.41%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.692'
$ws.Range('E38').Value = '  +0.24%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01962'
$ws.Range('E39').Value = '  -0.35%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.749'
$ws.Range('E40').Value = '  -1.43%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.4523'
$ws.Range('E41').Value = '  +0.32%  '

# Row 42
$ws.Range('D42').Value = '1.112.35'
$ws.Range('E42').Value = '  +0.11%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '6.073'
$ws.Range('E43').Value = '  -0.51%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '72.68'
$ws.Range('E44').Value = '  -2.04%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.8567'
$ws.Range('E45').Value = '  +0.31%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.000'
$ws.Range('E46').Value = '  +0.01%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '103.44'
$ws.Range('E47').Value = '  -0.01%  '

# Row 48
$ws.Range('B48').Value = 'SynthetixNetwork'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.120'
$ws.Range('E48').Value = '  +3.80%  '

# Row 49
$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.635'
$ws.Range('E49').Value = '  +0.28%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.841'
$ws.Range('E50').Value = '  -2.05%  '

# Row 51
$ws.Range('D51').Value = '2.021.87'
$ws.Range('E51').Value = '  -8.38%  '
